$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 156.656447
$ws.Range("H2").Value = 469.969341
$ws.Range("I2").Value = 0.0671576211124673
$ws.Range("J2").Value = 0.0671576211124673
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 85.89497033333333
$ws.Range("N2").Value = 257.684911
$ws.Range("O2").Value = 0.7848889718219874
$ws.Range("P2").Value = 0.7848889718219874
$ws.Range("Q2").Value = 13456.0008675904
$ws.Range("R2").Value = 121104.0078083136
$ws.Range("S2").Value = 0.05271127618497505
$ws.Range("T2").Value = 0.05271127618497505

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 156.656447
$ws.Range("H3").Value = 469.969341
$ws.Range("I3").Value = 0.0671576211124673
$ws.Range("J3").Value = 0.0671576211124673
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.726772333333333
$ws.Range("N3").Value = 5.180317
$ws.Range("O3").Value = 0.0157788582500353
$ws.Range("P3").Value = 0.0157788582500353
$ws.Range("Q3").Value = 270.5100185178996
$ws.Range("R3").Value = 2434.590166661097
$ws.Range("S3").Value = 0.001059670583943199
$ws.Range("T3").Value = 0.001059670583943199

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 156.656447
$ws.Range("H4").Value = 469.969341
$ws.Range("I4").Value = 0.0671576211124673
$ws.Range("J4").Value = 0.0671576211124673
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.21376166666667
$ws.Range("N4").Value = 54.641285
$ws.Range("O4").Value = 0.1664332685846793
$ws.Range("P4").Value = 0.1664332685846793
$ws.Range("Q4").Value = 2853.303189204798
$ws.Range("R4").Value = 25679.72870284319
$ws.Range("S4").Value = 0.0111772623921194
$ws.Range("T4").Value = 0.0111772623921194

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 156.656447
$ws.Range("H5").Value = 469.969341
$ws.Range("I5").Value = 0.0671576211124673
$ws.Range("J5").Value = 0.0671576211124673
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.600318333333334
$ws.Range("N5").Value = 10.800955
$ws.Range("O5").Value = 0.03289890134329811
$ws.Range("P5").Value = 0.03289890134329811
$ws.Range("Q5").Value = 564.0130781689617
$ws.Range("R5").Value = 5076.117703520655
$ws.Range("S5").Value = 0.002209411951429656
$ws.Range("T5").Value = 0.002209411951429656

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 42.300692
$ws.Range("H6").Value = 126.902076
$ws.Range("I6").Value = 0.01813403725498241
$ws.Range("J6").Value = 0.01813403725498241
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 85.89497033333333
$ws.Range("N6").Value = 257.684911
$ws.Range("O6").Value = 0.7848889718219874
$ws.Range("P6").Value = 0.7848889718219874
$ws.Range("Q6").Value = 3633.416684419471
$ws.Range("R6").Value = 32700.75015977523
$ws.Range("S6").Value = 0.01423320585604475
$ws.Range("T6").Value = 0.01423320585604475

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 42.300692
$ws.Range("H7").Value = 126.902076
$ws.Range("I7").Value = 0.01813403725498241
$ws.Range("J7").Value = 0.01813403725498241
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.726772333333333
$ws.Range("N7").Value = 5.180317
$ws.Range("O7").Value = 0.0157788582500353
$ws.Range("P7").Value = 0.0157788582500353
$ws.Range("Q7").Value = 73.04366462645466
$ws.Range("R7").Value = 657.3929816380919
$ws.Range("S7").Value = 0.0002861344033472266
$ws.Range("T7").Value = 0.0002861344033472266

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 42.300692
$ws.Range("H8").Value = 126.902076
$ws.Range("I8").Value = 0.01813403725498241
$ws.Range("J8").Value = 0.01813403725498241
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.21376166666667
$ws.Range("N8").Value = 54.641285
$ws.Range("O8").Value = 0.1664332685846793
$ws.Range("P8").Value = 0.1664332685846793
$ws.Range("Q8").Value = 770.4547224230733
$ws.Range("R8").Value = 6934.09250180766
$ws.Range("S8").Value = 0.003018107092983067
$ws.Range("T8").Value = 0.003018107092983067

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42.300692
$ws.Range("H9").Value = 126.902076
$ws.Range("I9").Value = 0.01813403725498241
$ws.Range("J9").Value = 0.01813403725498241
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.600318333333334
$ws.Range("N9").Value = 10.800955
$ws.Range("O9").Value = 0.03289890134329811
$ws.Range("P9").Value = 0.03289890134329811
$ws.Range("Q9").Value = 152.2959569202867
$ws.Range("R9").Value = 1370.66361228258
$ws.Range("S9").Value = 0.0005965899026073586
$ws.Range("T9").Value = 0.0005965899026073586

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2110.189616
$ws.Range("H10").Value = 6330.568848
$ws.Range("I10").Value = 0.9046248489651427
$ws.Range("J10").Value = 0.9046248489651426
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 85.89497033333333
$ws.Range("N10").Value = 257.684911
$ws.Range("O10").Value = 0.7848889718219874
$ws.Range("P10").Value = 0.7848889718219874
$ws.Range("Q10").Value = 181254.6744640281
$ws.Range("R10").Value = 1631292.070176252
$ws.Range("S10").Value = 0.7100300675888714
$ws.Range("T10").Value = 0.7100300675888714

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2110.189616
$ws.Range("H11").Value = 6330.568848
$ws.Range("I11").Value = 0.9046248489651427
$ws.Range("J11").Value = 0.9046248489651426
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.726772333333333
$ws.Range("N11").Value = 5.180317
$ws.Range("O11").Value = 0.0157788582500353
$ws.Range("P11").Value = 0.0157788582500353
$ws.Range("Q11").Value = 3643.817046996091
$ws.Range("R11").Value = 32794.35342296481
$ws.Range("S11").Value = 0.01427394726128058
$ws.Range("T11").Value = 0.01427394726128058

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2110.189616
$ws.Range("H12").Value = 6330.568848
$ws.Range("I12").Value = 0.9046248489651427
$ws.Range("J12").Value = 0.9046248489651426
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.21376166666667
$ws.Range("N12").Value = 54.641285
$ws.Range("O12").Value = 0.1664332685846793
$ws.Range("P12").Value = 0.1664332685846793
$ws.Range("Q12").Value = 38434.49073729885
$ws.Range("R12").Value = 345910.4166356897
$ws.Range("S12").Value = 0.1505596704561905
$ws.Range("T12").Value = 0.1505596704561905

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2110.189616
$ws.Range("H13").Value = 6330.568848
$ws.Range("I13").Value = 0.9046248489651427
$ws.Range("J13").Value = 0.9046248489651426
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.600318333333334
$ws.Range("N13").Value = 10.800955
$ws.Range("O13").Value = 0.03289890134329811
$ws.Range("P13").Value = 0.03289890134329811
$ws.Range("Q13").Value = 7597.354361294428
$ws.Range("R13").Value = 68376.18925164983
$ws.Range("S13").Value = 0.02976116365880018
$ws.Range("T13").Value = 0.02976116365880018

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 23.52144266666667
$ws.Range("H14").Value = 70.564328
$ws.Range("I14").Value = 0.01008349266740757
$ws.Range("J14").Value = 0.01008349266740757
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 85.89497033333333
$ws.Range("N14").Value = 257.684911
$ws.Range("O14").Value = 0.7848889718219874
$ws.Range("P14").Value = 0.7848889718219874
$ws.Range("Q14").Value = 2020.373620050534
$ws.Range("R14").Value = 18183.36258045481
$ws.Range("S14").Value = 0.007914422192096078
$ws.Range("T14").Value = 0.007914422192096076

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 23.52144266666667
$ws.Range("H15").Value = 70.564328
$ws.Range("I15").Value = 0.01008349266740757
$ws.Range("J15").Value = 0.01008349266740757
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.726772333333333
$ws.Range("N15").Value = 5.180317
$ws.Range("O15").Value = 0.0157788582500353
$ws.Range("P15").Value = 0.0157788582500353
$ws.Range("Q15").Value = 40.61617643688622
$ws.Range("R15").Value = 365.545587931976
$ws.Range("S15").Value = 0.0001591060014642944
$ws.Range("T15").Value = 0.0001591060014642944

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 23.52144266666667
$ws.Range("H16").Value = 70.564328
$ws.Range("I16").Value = 0.01008349266740757
$ws.Range("J16").Value = 0.01008349266740757
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 18.21376166666667
$ws.Range("N16").Value = 54.641285
$ws.Range("O16").Value = 0.1664332685846793
$ws.Range("P16").Value = 0.1664332685846793
$ws.Range("Q16").Value = 428.4139507868312
$ws.Range("R16").Value = 3855.72555708148
$ws.Range("S16").Value = 0.001678228643386288
$ws.Range("T16").Value = 0.001678228643386288

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 23.52144266666667
$ws.Range("H17").Value = 70.564328
$ws.Range("I17").Value = 0.01008349266740757
$ws.Range("J17").Value = 0.01008349266740757
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.600318333333334
$ws.Range("N17").Value = 10.800955
$ws.Range("O17").Value = 0.03289890134329811
$ws.Range("P17").Value = 0.03289890134329811
$ws.Range("Q17").Value = 84.6846812592489
$ws.Range("R17").Value = 762.16213133324
$ws.Range("S17").Value = 0.0003317358304609116
$ws.Range("T17").Value = 0.0003317358304609115
